# Update "想去人数" (interest count) figures in the "F" column across the
# sheets of the 苏州-漫展信息 workbook, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 1044
$wsExhibit.Range("F6").Value  = 166
$wsExhibit.Range("F8").Value  = 200
$wsExhibit.Range("F9").Value  = 378
$wsExhibit.Range("F12").Value = 527
$wsExhibit.Range("F13").Value = 152
$wsExhibit.Range("F14").Value = 12366
$wsExhibit.Range("F15").Value = 98
$wsExhibit.Range("F16").Value = 5477

# --- Sheet "演出" ---------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 115

# --- Sheet "全部类型" ------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 115
$wsAll.Range("F7").Value  = 1044
$wsAll.Range("F8").Value  = 166
$wsAll.Range("F10").Value = 200
$wsAll.Range("F11").Value = 378
$wsAll.Range("F14").Value = 527
$wsAll.Range("F15").Value = 152
$wsAll.Range("F16").Value = 12366
$wsAll.Range("F18").Value = 98
$wsAll.Range("F19").Value = 5477

$wb.Save()
